$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$changes = @(
    @("H106", 3917.6667),
    @("I106", 2625),
    @("K106", 2625),
    @("M106", -1994),
    @("H129", 905.77026),
    @("J129", 884.1739),
    @("L129", 2652.5217),
    @("N129", -12652.5217),
    @("H137", 1669.4445),
    @("I137", 1114.5385),
    @("J137", 2184.7144),
    @("K137", 3343.6155),
    @("L137", 6554.1432),
    @("M137", -793.6155000000003),
    @("N137", -11654.1432),
    @("H138", 2842.16),
    @("I138", 2954.9524),
    @("J138", 2250),
    @("K138", 8864.8572),
    @("L138", 6750),
    @("M138", -3724.8572),
    @("N138", -17030)
)
foreach ($chg in $changes) {
    $ref = $chg[0]
    $val = $chg[1]
    if ($val -eq $null) {
        $ws.Range($ref).ClearContents()
    } else {
        $ws.Range($ref).Value = $val
    }
}

$ws = $wb.Worksheets.Item("ARM")
$changes = @(
    @("H32", 3591.4092),
    @("I32", 2341.796),
    @("J32", 7193.2354),
    @("K32", 2341.796),
    @("L32", 7193.2354),
    @("M32", -2054.796),
    @("N32", -7767.2354),
    @("H45", 1530.5),
    @("I45", 987.75),
    @("K45", 987.75),
    @("M45", -610.75),
    @("H74", 972.6667),
    @("J74", 1156.5),
    @("L74", 1156.5),
    @("N74", -2904.5),
    @("H77", 972.6667),
    @("J77", 1156.5),
    @("L77", 5782.5),
    @("N77", -14518.5),
    @("H104", 39999.75),
    @("J104", 39999.75),
    @("L104", 39999.75),
    @("N104", -46987.75),
    @("H110", 1499.7273),
    @("I110", 1117.9286),
    @("J110", 2167.875),
    @("K110", 1117.9286),
    @("L110", 2167.875),
    @("M110", 927.0714),
    @("N110", -6257.875)
)
foreach ($chg in $changes) {
    $ref = $chg[0]
    $val = $chg[1]
    if ($val -eq $null) {
        $ws.Range($ref).ClearContents()
    } else {
        $ws.Range($ref).Value = $val
    }
}

$ws = $wb.Worksheets.Item("BSM")
$changes = @(
    @("H122", 0),
    @("J122", 0),
    @("L122", 0),
    @("N122", $null),
    @("H137", 0),
    @("J137", 0),
    @("L137", 0),
    @("N137", $null)
)
foreach ($chg in $changes) {
    $ref = $chg[0]
    $val = $chg[1]
    if ($val -eq $null) {
        $ws.Range($ref).ClearContents()
    } else {
        $ws.Range($ref).Value = $val
    }
}

$ws = $wb.Worksheets.Item("CRP")
$changes = @(
    @("H31", 2520),
    @("I31", 1198.4546),
    @("J31", 6673.4287),
    @("K31", 1198.4546),
    @("L31", 6673.4287),
    @("M31", -903.4546),
    @("N31", -7263.4287),
    @("H34", 2520),
    @("I34", 1198.4546),
    @("J34", 6673.4287),
    @("K34", 1198.4546),
    @("L34", 6673.4287),
    @("M34", -996.4546),
    @("N34", -7077.4287),
    @("H134", 962.9091),
    @("I134", 859.2),
    @("K134", 2577.6),
    @("M134", -42.60000000000036)
)
foreach ($chg in $changes) {
    $ref = $chg[0]
    $val = $chg[1]
    if ($val -eq $null) {
        $ws.Range($ref).ClearContents()
    } else {
        $ws.Range($ref).Value = $val
    }
}

$ws = $wb.Worksheets.Item("CUL")
$changes = @(
    @("H5", 745.44446),
    @("I5", 634),
    @("K5", 1902),
    @("M5", -1790),
    @("H22", 2000.5),
    @("I22", 2000.5),
    @("K22", 6001.5),
    @("M22", -5832.5),
    @("H27", 2000.5),
    @("I27", 2000.5),
    @("K27", 6001.5),
    @("M27", -5899.5),
    @("H37", 99932.664),
    @("J37", 99932.664),
    @("L37", 299797.992),
    @("N37", -300021.992),
    @("H43", 0),
    @("J43", 0),
    @("L43", 0),
    @("N43", $null),
    @("H45", 983.6),
    @("J45", 983.6),
    @("L45", 2950.8),
    @("N45", -4014.8),
    @("H57", 3000),
    @("J57", 4000),
    @("L57", 12000),
    @("N57", -13118),
    @("H62", 3500),
    @("I62", 3000),
    @("K62", 9000),
    @("M62", -8314),
    @("H65", 3500),
    @("I65", 3000),
    @("K65", 27000),
    @("M65", -23568),
    @("H113", 8599.691999999999),
    @("J113", 1019.2),
    @("L113", 3057.6),
    @("N113", -7397.6),
    @("H116", 3000),
    @("J116", 3000),
    @("L116", 9000),
    @("N116", -15884),
    @("H131", 770.49),
    @("I131", 332.33334),
    @("J131", 813.82416),
    @("K131", 997.0000200000001),
    @("L131", 2441.47248),
    @("M131", 4042.99998),
    @("N131", -12521.47248),
    @("H135", 745.44446),
    @("I135", 634),
    @("K135", 5706),
    @("M135", -3171),
    @("H140", 1792.5172),
    @("I140", 932.25),
    @("K140", 2796.75),
    @("M140", 2383.25)
)
foreach ($chg in $changes) {
    $ref = $chg[0]
    $val = $chg[1]
    if ($val -eq $null) {
        $ws.Range($ref).ClearContents()
    } else {
        $ws.Range($ref).Value = $val
    }
}

$ws = $wb.Worksheets.Item("GSM")
$changes = @(
    @("H70", 5180.8),
    @("I70", 5833.1665),
    @("K70", 5833.1665),
    @("M70", -5563.1665),
    @("H73", 5180.8),
    @("I73", 5833.1665),
    @("K73", 5833.1665),
    @("M73", -4897.1665),
    @("H132", 5346.25),
    @("I132", 4317.5625),
    @("J132", 9461),
    @("K132", 12952.6875),
    @("L132", 28383),
    @("M132", -10422.6875),
    @("N132", -33443)
)
foreach ($chg in $changes) {
    $ref = $chg[0]
    $val = $chg[1]
    if ($val -eq $null) {
        $ws.Range($ref).ClearContents()
    } else {
        $ws.Range($ref).Value = $val
    }
}

$ws = $wb.Worksheets.Item("LTW")
$changes = @(
    @("H40", 11599.8),
    @("I40", 6334.6665),
    @("K40", 6334.6665),
    @("M40", -6198.6665),
    @("H46", 2750),
    @("J46", 2750),
    @("L46", 2750),
    @("N46", -3126),
    @("H132", 1838.5834),
    @("I132", 1502.6666),
    @("J132", 1950.5555),
    @("K132", 4507.9998),
    @("L132", 5851.666499999999),
    @("M132", -1977.9998),
    @("N132", -10911.6665)
)
foreach ($chg in $changes) {
    $ref = $chg[0]
    $val = $chg[1]
    if ($val -eq $null) {
        $ws.Range($ref).ClearContents()
    } else {
        $ws.Range($ref).Value = $val
    }
}
